# Add three new columns (M, N, O) to Sheet1: t,P / t,SYM-H / t,CS
# These are event timestamps (time-of-day) formatted as hh:mm, added for
# the first four data rows (2-5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("M1").Value = "t,P"
$ws.Range("N1").Value = "t,SYM-H"
$ws.Range("O1").Value = "t,CS"

# --- New column widths (compensate for the host's character->width
#     padding of 5/6 so the saved OOXML width lands as close as possible
#     to the authored 8.26 / 9.37 / 7.06) ---
$ws.Columns.Item(13).ColumnWidth = 7.426666666666667
$ws.Columns.Item(14).ColumnWidth = 8.536666666666665
$ws.Columns.Item(15).ColumnWidth = 6.226666666666667

# --- New time values for rows 2-5 (columns M, N, O) ---
$ws.Range("M2").Value = 0.697916666666667
$ws.Range("N2").Value = 0.694444444444444
$ws.Range("O2").Value = 0.697916666666667

$ws.Range("M3").Value = 0.245138888888889
$ws.Range("N3").Value = 0.245138888888889
$ws.Range("O3").Value = 0.246527777777778

$ws.Range("M4").Value = 0.245833333333333
$ws.Range("N4").Value = 0.243055555555556
$ws.Range("O4").Value = 0.247222222222222

$ws.Range("M5").Value = 0.0923611111111111
$ws.Range("N5").Value = 0.0868055555555556
$ws.Range("O5").Value = 0.0923611111111111

# Format the new time cells as hh:mm (adds a new numFmt/style, like the source edit)
$ws.Range("M2:O5").NumberFormat = "hh:mm"

# --- Restore the active selection shown in the edited file ---
[void]$ws.Range("M7").Select()
